$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 32 and 33 (columns B:AC)
$ws.Range("B32").Value = 5400042
$ws.Range("B33").Value = 5400043
$ws.Range("F32").Value = "Melbourne City"
$ws.Range("F33").Value = "Sydney FC"
$ws.Range("G32").Value = "Macarthur FC"
$ws.Range("G33").Value = "Central Coast Mariners"
$ws.Range("H32").Value = 6
$ws.Range("H33").Value = 3
$ws.Range("I32").Value = 1
$ws.Range("I33").Value = 2
$ws.Range("K32").Value = 1.533
$ws.Range("K33").Value = 2.15
$ws.Range("L32").Value = 4.2
$ws.Range("L33").Value = 3.4
$ws.Range("M32").Value = 6
$ws.Range("M33").Value = 3.4
$ws.Range("N32").Value = 1.333
$ws.Range("N33").Value = 2.2
$ws.Range("O32").Value = 5
$ws.Range("O33").Value = 3.75
$ws.Range("P32").Value = 9.5
$ws.Range("P33").Value = 3
$ws.Range("Q32").Value = -1.5
$ws.Range("Q33").Value = -0.25
$ws.Range("R32").Value = 1.85
$ws.Range("R33").Value = 1.975
$ws.Range("S32").Value = 2
$ws.Range("S33").Value = 1.875
$ws.Range("T32").Value = 3.25
$ws.Range("T33").Value = 3
$ws.Range("U32").Value = 1.875
$ws.Range("U33").Value = 1.95
$ws.Range("V32").Value = 1.975
$ws.Range("V33").Value = 1.9
$ws.Range("W32").Value = 0.333
$ws.Range("W33").Value = 1.2
$ws.Range("Z32").Value = 0.8500000000000001
$ws.Range("Z33").Value = 0.9750000000000001
$ws.Range("AB32").Value = 0.875
$ws.Range("AB33").Value = 0.95

# Swap rows 50 and 51 (columns B:AC)
$ws.Range("B50").Value = 5404699
$ws.Range("B51").Value = 5400048
$ws.Range("F50").Value = "Newcastle Jets"
$ws.Range("F51").Value = "Melbourne City"
$ws.Range("G50").Value = "Western United FC"
$ws.Range("G51").Value = "Sydney FC"
$ws.Range("H50").Value = 1
$ws.Range("H51").Value = 3
$ws.Range("I50").Value = 3
$ws.Range("I51").Value = 2
$ws.Range("J50").Value = "A"
$ws.Range("J51").Value = "H"
$ws.Range("K50").Value = 2.5
$ws.Range("K51").Value = 1.7
$ws.Range("L50").Value = 3.1
$ws.Range("L51").Value = 4
$ws.Range("M50").Value = 2.8
$ws.Range("M51").Value = 4
$ws.Range("N50").Value = 2.1
$ws.Range("N51").Value = 1.8
$ws.Range("O50").Value = 3.6
$ws.Range("O51").Value = 3.8
$ws.Range("P50").Value = 3.3
$ws.Range("P51").Value = 4
$ws.Range("Q50").Value = -0.25
$ws.Range("Q51").Value = -0.75
$ws.Range("R50").Value = 1.94
$ws.Range("R51").Value = 2.025
$ws.Range("S50").Value = 1.96
$ws.Range("S51").Value = 1.825
$ws.Range("U50").Value = 1.975
$ws.Range("U51").Value = 1.85
$ws.Range("V50").Value = 1.875
$ws.Range("V51").Value = 2
$ws.Range("W50").Value = -1
$ws.Range("W51").Value = 0.8
$ws.Range("Y50").Value = 2.3
$ws.Range("Y51").Value = -1
$ws.Range("Z50").Value = -1
$ws.Range("Z51").Value = 0.5125
$ws.Range("AA50").Value = 0.96
$ws.Range("AA51").Value = -0.5
$ws.Range("AB50").Value = 0.9750000000000001
$ws.Range("AB51").Value = 0.8500000000000001

# Swap rows 57 and 58 (columns B:AC)
$ws.Range("B57").Value = 5404704
$ws.Range("B58").Value = 5404706
$ws.Range("F57").Value = "Western Sydney Wanderers"
$ws.Range("F58").Value = "Western United FC"
$ws.Range("G57").Value = "Central Coast Mariners"
$ws.Range("G58").Value = "Perth Glory"
$ws.Range("I57").Value = 0
$ws.Range("I58").Value = 1
$ws.Range("K57").Value = 2.4
$ws.Range("K58").Value = 2.1
$ws.Range("L57").Value = 3.5
$ws.Range("L58").Value = 3.4
$ws.Range("M57").Value = 2.6
$ws.Range("M58").Value = 3.2
$ws.Range("N57").Value = 2.2
$ws.Range("N58").Value = 2.1
$ws.Range("O57").Value = 3.75
$ws.Range("O58").Value = 3.6
$ws.Range("P57").Value = 3.1
$ws.Range("P58").Value = 3.3
$ws.Range("R57").Value = 1.875
$ws.Range("R58").Value = 1.85
$ws.Range("S57").Value = 1.975
$ws.Range("S58").Value = 2
$ws.Range("U57").Value = 1.875
$ws.Range("U58").Value = 1.85
$ws.Range("V57").Value = 1.975
$ws.Range("V58").Value = 2
$ws.Range("W57").Value = 1.2
$ws.Range("W58").Value = 1.1
$ws.Range("Z57").Value = 0.875
$ws.Range("Z58").Value = 0.8500000000000001
$ws.Range("AB57").Value = -1
$ws.Range("AB58").Value = 0.425
$ws.Range("AC57").Value = 0.9750000000000001
$ws.Range("AC58").Value = -0.5

# Swap rows 68 and 69 (columns B:AC)
$ws.Range("B68").Value = 5404713
$ws.Range("B69").Value = 5404714
$ws.Range("F68").Value = "Brisbane Roar"
$ws.Range("F69").Value = "Newcastle Jets"
$ws.Range("G68").Value = "Western United FC"
$ws.Range("G69").Value = "Perth Glory"
$ws.Range("H68").Value = 1
$ws.Range("H69").Value = 2
$ws.Range("I68").Value = 0
$ws.Range("I69").Value = 2
$ws.Range("J68").Value = "H"
$ws.Range("J69").Value = "D"
$ws.Range("L68").Value = 3.5
$ws.Range("L69").Value = 3.3
$ws.Range("M68").Value = 2.8
$ws.Range("M69").Value = 3
$ws.Range("N68").Value = 2.3
$ws.Range("N69").Value = 2.15
$ws.Range("O68").Value = 3.6
$ws.Range("O69").Value = 3.75
$ws.Range("P68").Value = 2.9
$ws.Range("P69").Value = 3.2
$ws.Range("R68").Value = 2.025
$ws.Range("R69").Value = 1.9
$ws.Range("S68").Value = 1.825
$ws.Range("S69").Value = 1.95
$ws.Range("U68").Value = 2.025
$ws.Range("U69").Value = 1.9
$ws.Range("V68").Value = 1.825
$ws.Range("V69").Value = 1.95
$ws.Range("W68").Value = 1.3
$ws.Range("W69").Value = -1
$ws.Range("X68").Value = -1
$ws.Range("X69").Value = 2.75
$ws.Range("Z68").Value = 1.025
$ws.Range("Z69").Value = -0.5
$ws.Range("AA68").Value = -1
$ws.Range("AA69").Value = 0.475
$ws.Range("AB68").Value = -1
$ws.Range("AB69").Value = 0.8999999999999999
$ws.Range("AC68").Value = 0.825
$ws.Range("AC69").Value = -1

# Swap rows 180 and 181 (columns B:AC)
$ws.Range("B180").Value = 7646750
$ws.Range("B181").Value = 7646749
$ws.Range("F180").Value = "Perth Glory"
$ws.Range("F181").Value = "Brisbane Roar"
$ws.Range("G180").Value = "Wellington Phoenix"
$ws.Range("G181").Value = "Newcastle Jets"
$ws.Range("I180").Value = 4
$ws.Range("I181").Value = 2
$ws.Range("J180").Value = "A"
$ws.Range("J181").Value = "H"
$ws.Range("K180").Value = 2.45
$ws.Range("K181").Value = 1.909
$ws.Range("L180").Value = 3.75
$ws.Range("L181").Value = 4
$ws.Range("M180").Value = 2.55
$ws.Range("M181").Value = 3.4
$ws.Range("N180").Value = 3.1
$ws.Range("N181").Value = 2.4
$ws.Range("O180").Value = 3.8
$ws.Range("O181").Value = 4
$ws.Range("P180").Value = 2.05
$ws.Range("P181").Value = 2.6
$ws.Range("Q180").Value = 0.25
$ws.Range("Q181").Value = 0
$ws.Range("R180").Value = 2
$ws.Range("R181").Value = 1.83
$ws.Range("S180").Value = 1.85
$ws.Range("S181").Value = 2.07
$ws.Range("T180").Value = 3
$ws.Range("T181").Value = 3.25
$ws.Range("U180").Value = 1.925
$ws.Range("U181").Value = 1.9
$ws.Range("V180").Value = 1.925
$ws.Range("V181").Value = 1.95
$ws.Range("W180").Value = -1
$ws.Range("W181").Value = 1.4
$ws.Range("Y180").Value = 1.05
$ws.Range("Y181").Value = -1
$ws.Range("Z180").Value = -1
$ws.Range("Z181").Value = 0.8300000000000001
$ws.Range("AA180").Value = 0.8500000000000001
$ws.Range("AA181").Value = -1
$ws.Range("AB180").Value = 0.925
$ws.Range("AB181").Value = 0.8999999999999999
# Individual cell updates (rows 214, 216, 217, 218, 219)
$ws.Range("N214").Value = 1.909
$ws.Range("O214").Value = 4
$ws.Range("P214").Value = 3.6
$ws.Range("R214").Value = 1.95
$ws.Range("S214").Value = 1.95
$ws.Range("U214").Value = 2.025
$ws.Range("V214").Value = 1.825

$ws.Range("R216").Value = 1.83
$ws.Range("S216").Value = 2.07

$ws.Range("N217").Value = 2.55
$ws.Range("P217").Value = 2.6
$ws.Range("R217").Value = 1.88
$ws.Range("S217").Value = 2.02

$ws.Range("O218").Value = 4
$ws.Range("Q218").Value = -0.25
$ws.Range("R218").Value = 1.83
$ws.Range("S218").Value = 2.07

$ws.Range("R219").Value = 1.99
$ws.Range("S219").Value = 1.91

Write-Output "Applied Australia ALeague updates"
